$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $A, $B, $C, $D, $E, $F, $G, $H, $K) {
    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H

    if ($K -ne $null) {
        # K column stores decimal-looking values as text, so force text
        # to avoid Excel auto-converting them into numbers.
        $ws.Cells.Item($Row, 11).NumberFormat = "@"
        $ws.Cells.Item($Row, 11).Value = $K
    }
}

Set-Row 116 "Anomaly no noise" "yxuau6po" "Training phase" 1 `
    "['Purple', 'Green']" `
    "[['Red', 'Blue'], ['Blue', 'Yellow']]" `
    "[None, None]" `
    "['8', '8']" `
    "0.16"

Set-Row 117 "Transmission correct" "2nj8y4ca" "Training phase" 1 `
    "['Purple', 'Green']" `
    "[['Red', ''], ['Blue', '']]" `
    "[None, None]" `
    "['2', '5']" `
    "0.07"

Set-Row 118 "Transmission correct" "txpcyowh" "Training phase" 1 `
    "['Purple', 'Green']" `
    "[['Red', ''], ['Blue', '']]" `
    "[None, None]" `
    "['2', '5']" `
    $null

Set-Row 119 "Transmission correct" "txpcyowh" "Training phase" 2 `
    "['Green', 'Purple']" `
    "[['Red', ''], ['Blue', '']]" `
    "[None, None]" `
    "['2', '5']" `
    $null

Set-Row 120 "Transmission correct" "txpcyowh" "Training phase" 3 `
    "['Purple', 'Green', 'Green']" `
    "[['Red', ''], ['Blue', ''], ['Blue', '']]" `
    "[None, None, None]" `
    "['2', '5', '5']" `
    $null

Set-Row 121 "Transmission correct" "txpcyowh" "Training phase" 4 `
    "['Green', 'Purple', 'Purple']" `
    "[['Blue', ''], ['Blue', ''], ['Blue', '']]" `
    "[None, None, None]" `
    "['5', '5', '5']" `
    $null

Set-Row 122 "Transmission correct" "txpcyowh" "Test 1" 1 `
    "['Green', 'Yellow', 'Purple', 'Red', 'Orange', 'Blue']" `
    "[['Red', ''], ['Red', ''], ['Blue', ''], ['Red', ''], ['Red', ''], ['Red', '']]" `
    "[None, None, None, None, None, None]" `
    "['0', '0', '0', '0', '0', '0']" `
    $null

Set-Row 123 "Transmission correct" "txpcyowh" "Exploration phase" 1 `
    "['Orange', 'Purple']" `
    "[['Red', ''], ['Red', '']]" `
    "[None, None]" `
    "['2', '2']" `
    $null

Set-Row 124 "Transmission correct" "txpcyowh" "Exploration phase" 2 `
    "['Orange', 'Green']" `
    "[['Red', ''], ['Red', '']]" `
    "[None, None]" `
    "['2', '2']" `
    $null

Set-Row 125 "Transmission correct" "txpcyowh" "Exploration phase" 3 `
    "['Green', 'Red', 'Green']" `
    "[['Red', ''], ['Red', ''], ['Red', '']]" `
    "[None, None, None]" `
    "['2', '2', '2']" `
    $null

Set-Row 126 "Transmission correct" "txpcyowh" "Exploration phase" 4 `
    "['Blue', 'Purple', 'Purple']" `
    "[['Red', ''], ['Red', ''], ['Red', '']]" `
    "[None, None, None]" `
    "['2', '2', '2']" `
    $null

Set-Row 127 "Transmission correct" "txpcyowh" "Exploration phase" 5 `
    "['Purple', 'Green', 'Yellow']" `
    "[['Red', ''], ['Red', ''], ['Red', '']]" `
    "[None, None, None]" `
    "['2', '2', '2']" `
    $null

Set-Row 128 "Transmission correct" "txpcyowh" "Test 2" 1 `
    "['Green', 'Yellow', 'Purple', 'Red', 'Orange', 'Blue']" `
    "[['Red', ''], ['Red', ''], ['Red', ''], ['Red', ''], ['Red', ''], ['Red', '']]" `
    "[None, None, None, None, None, None]" `
    "['0', '0', '0', '0', '0', '0']" `
    "0.67"
